$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 3500
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 56
$ws.Range("H56").Value = 3500
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
# Row 125
$ws.Range("H125").Value = 5412.222
$ws.Range("J125").Value = 6500
$ws.Range("L125").Value = 58500
$ws.Range("N125").Value = -63420
# Row 132
$ws.Range("H132").Value = 2503.4783
$ws.Range("I132").Value = 1789.7142
$ws.Range("J132").Value = 9998
$ws.Range("K132").Value = 5369.142599999999
$ws.Range("L132").Value = 29994
$ws.Range("M132").Value = -2839.142599999999
$ws.Range("N132").Value = -35054
# Row 137
$ws.Range("H137").Value = 441499.2
$ws.Range("I137").Value = 2121
$ws.Range("J137").Value = 807647.7
$ws.Range("K137").Value = 6363
$ws.Range("L137").Value = 2422943.1
$ws.Range("M137").Value = -3813
$ws.Range("N137").Value = -2428043.1
# Row 138
$ws.Range("H138").Value = 1672.9333
$ws.Range("I138").Value = 1468.6666
$ws.Range("J138").Value = 2490
$ws.Range("K138").Value = 4405.9998
$ws.Range("L138").Value = 7470
$ws.Range("M138").Value = 734.0002000000004
$ws.Range("N138").Value = -17750

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6401.15
$ws.Range("I32").Value = 3445.7317
$ws.Range("K32").Value = 3445.7317
$ws.Range("M32").Value = -3158.7317
# Row 43
$ws.Range("H43").Value = 25213.5
$ws.Range("I43").Value = 25050
$ws.Range("J43").Value = 25377
$ws.Range("K43").Value = 25050
$ws.Range("L43").Value = 25377
$ws.Range("M43").Value = -24737
$ws.Range("N43").Value = -26003
# Row 74
$ws.Range("H74").Value = 1464.625
$ws.Range("I74").Value = 1431.0476
$ws.Range("K74").Value = 1431.0476
$ws.Range("M74").Value = -557.0476000000001
# Row 77
$ws.Range("H77").Value = 1464.625
$ws.Range("I77").Value = 1431.0476
$ws.Range("K77").Value = 7155.238
$ws.Range("M77").Value = -2787.238
# Row 104
$ws.Range("H104").Value = 48945.332
$ws.Range("J104").Value = 48945.332
$ws.Range("L104").Value = 48945.332
$ws.Range("N104").Value = -55933.332
# Row 109
$ws.Range("H109").Value = 54999.668
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57774
# Row 132
$ws.Range("H132").Value = 1661.5555
$ws.Range("I132").Value = 1259.8334
$ws.Range("J132").Value = 2465
$ws.Range("K132").Value = 3779.5002
$ws.Range("L132").Value = 7395
$ws.Range("M132").Value = -1249.5002
$ws.Range("N132").Value = -12455

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2681.7778
$ws.Range("I94").Value = 2814.8096
$ws.Range("K94").Value = 2814.8096
$ws.Range("M94").Value = -2363.8096
# Row 108
$ws.Range("H108").Value = 92662.78
$ws.Range("J108").Value = 92662.78
$ws.Range("L108").Value = 92662.78
$ws.Range("N108").Value = -100342.78

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Range("H60").Value = 2861.0908
$ws.Range("I60").Value = 2861.0908
$ws.Range("K60").Value = 2861.0908
$ws.Range("M60").Value = -2350.0908
# Row 63
$ws.Range("H63").Value = 55000
$ws.Range("J63").Value = 55000
$ws.Range("L63").Value = 55000
$ws.Range("N63").Value = -56372
# Row 66
$ws.Range("H66").Value = 55000
$ws.Range("J66").Value = 55000
$ws.Range("L66").Value = 165000
$ws.Range("N66").Value = -171864

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 247.5
$ws.Range("I2").Value = 475
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 2850
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = -2737
$ws.Range("N2").Value = -346
# Row 56
$ws.Range("H56").Value = 6331.6
$ws.Range("I56").Value = 6331.6
$ws.Range("K56").Value = 6331.6
$ws.Range("M56").Value = -5801.6
# Row 86
$ws.Range("H86").Value = 9717
$ws.Range("I86").Value = 20694.4
$ws.Range("K86").Value = 62083.2
$ws.Range("M86").Value = -60897.2
# Row 89
$ws.Range("H89").Value = 9717
$ws.Range("I89").Value = 20694.4
$ws.Range("K89").Value = 186249.6
$ws.Range("M89").Value = -180321.6
# Row 103
$ws.Range("H103").Value = 1140.1
$ws.Range("I103").Value = 1236.1428
$ws.Range("J103").Value = 916
$ws.Range("K103").Value = 3708.4284
$ws.Range("L103").Value = 2748
$ws.Range("M103").Value = -2829.4284
$ws.Range("N103").Value = -4506
# Row 120
$ws.Range("H120").Value = 50001
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 50001
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 150003
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -159679

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 29330.5
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 132
$ws.Range("H132").Value = 3529.6538
$ws.Range("I132").Value = 2924.8823
$ws.Range("K132").Value = 8774.6469
$ws.Range("M132").Value = -6244.6469

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3126.0667
$ws.Range("I7").Value = 2344.5
$ws.Range("J7").Value = 3809.9375
$ws.Range("K7").Value = 2344.5
$ws.Range("L7").Value = 3809.9375
$ws.Range("M7").Value = -2232.5
$ws.Range("N7").Value = -4033.9375
# Row 58
$ws.Range("H58").Value = 2888
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 5000
$ws.Range("N58").Value = -5520
# Row 126
$ws.Range("H126").Value = 3126.0667
$ws.Range("I126").Value = 2344.5
$ws.Range("J126").Value = 3809.9375
$ws.Range("K126").Value = 7033.5
$ws.Range("L126").Value = 11429.8125
$ws.Range("M126").Value = -4563.5
$ws.Range("N126").Value = -16369.8125
# Row 132
$ws.Range("H132").Value = 1522.4231
$ws.Range("I132").Value = 1390.8695
$ws.Range("K132").Value = 4172.6085
$ws.Range("M132").Value = -1642.6085
# Row 136
$ws.Range("H136").Value = 4281.6665
$ws.Range("I136").Value = 4718.2
$ws.Range("J136").Value = 2099
$ws.Range("K136").Value = 14154.6
$ws.Range("L136").Value = 6297
$ws.Range("M136").Value = -11604.6
$ws.Range("N136").Value = -11397

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 19271
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 19271
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 19271
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -20769
# Row 72
$ws.Range("H72").Value = 19271
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 19271
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 57813
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -65301
# Row 75
$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 20000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -19064
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 20000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 60000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -55320
$ws.Range("N78").ClearContents()
# Row 132
$ws.Range("H132").Value = 1768.5
$ws.Range("I132").Value = 1413.5555
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 4240.666499999999
$ws.Range("L132").Value = 8499.999899999999
$ws.Range("M132").Value = -1710.666499999999
$ws.Range("N132").Value = -13559.9999
# Row 136
$ws.Range("H136").Value = 848.75
$ws.Range("I136").Value = 505.7143
$ws.Range("K136").Value = 1517.1429
$ws.Range("M136").Value = 1032.8571
